# Adds the "assessment_criteria" family of sheets (add/search/view/filter/edit),
# mirroring the existing "qualification_types" family, and tweaks a few
# existing qualification_types sheets (selection + one status value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Small edits on the existing qualification_types sheets
# ---------------------------------------------------------------------------

$wsAddQT = $wb.Worksheets.Item("add_new_qualification_types")
$wsAddQT.Range("E3").Value = "Inactive"

$wsSearchQT = $wb.Worksheets.Item("search_qualification_types")
$wsViewQT   = $wb.Worksheets.Item("view_qualification_types")
$wsFilterQT = $wb.Worksheets.Item("filter_qualification_types")
$wsEditQT   = $wb.Worksheets.Item("edit_qualification_types")

$wsSearchQT.Range("A1:B1").Select() | Out-Null
$wsViewQT.Range("A1:B3").Select() | Out-Null
$wsFilterQT.Range("A1:C2").Select() | Out-Null
$wsEditQT.Range("A2:A3").Select() | Out-Null

$wsAddQT.Range("G5").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Add the five new assessment_criteria sheets at the end of the workbook
# ---------------------------------------------------------------------------

function Add-SheetAtEnd($workbook, $name) {
    $last = $workbook.Worksheets.Item($workbook.Worksheets.Count)
    $newSheet = $workbook.Worksheets.Add($null, $last)
    $newSheet.Name = $name
    return $newSheet
}

$sAdd    = Add-SheetAtEnd $wb "add_new_assessment_criteria"
$sSearch = Add-SheetAtEnd $wb "search_assessment_criteria"
$sView   = Add-SheetAtEnd $wb "view_assessment_criteria"
$sFilter = Add-SheetAtEnd $wb "filter_assessment_criteria"
$sEdit   = Add-SheetAtEnd $wb "edit_assessment_criteria"

# --- add_new_assessment_criteria -------------------------------------------
$sAdd.Range("A2").Value = "Auto_AC_001"
$sAdd.Range("A3").Value = "Auto_AC_002"
$sAdd.Range("C1").Value = "sub assessment1"
$sAdd.Range("D1").Value = "sub assessment2"
$sAdd.Range("B2").Value = "Final Exam"
$sAdd.Range("C2").Value = "Theory"
$sAdd.Range("D2").Value = "Viva"
$sAdd.Range("B3").Value = "Mid Exam"
$sAdd.Range("A1").Value = "code"
$sAdd.Range("B1").Value = "name"
$sAdd.Range("E1").Value = "status"
$sAdd.Range("F1").Value = "runmode"
$sAdd.Range("E2").Value = "Active"
$sAdd.Range("F2").Value = "Y"
$sAdd.Range("C3").Value = "Theory"
$sAdd.Range("D3").Value = "Viva"
$sAdd.Range("E3").Value = "Inactive"
$sAdd.Range("F3").Value = "Y"

$sAdd.Columns.Item(1).ColumnWidth = 14.592447916666666
$sAdd.Columns.Item(2).ColumnWidth = 13.592447916666666
$sAdd.Columns.Item(3).ColumnWidth = 15.307291666666666
$sAdd.Columns.Item(4).ColumnWidth = 15.736979166666666

# --- edit_assessment_criteria (written before the search/filter sheets so
# the shared-string allocation order matches the source workbook) ----------
$sEdit.Range("B2").Value = "Auto_UPD_AC_001"
$sEdit.Range("B3").Value = "Auto_UPD_AC_002"
$sEdit.Range("C2").Value = "UPD_Final Exam"
$sEdit.Range("C3").Value = "UPD_Mid Exam"

# --- search_assessment_criteria --------------------------------------------
$sSearch.Range("A2").Value = "AC_001"
$sSearch.Range("A3").Value = "Final"
$sSearch.Range("A1").Value = "keyword"
$sSearch.Range("B1").Value = "runmode"
$sSearch.Range("B2").Value = "Y"
$sSearch.Range("B3").Value = "Y"

# --- filter_assessment_criteria ---------------------------------------------
$sFilter.Range("B2").Value = "Auto_AC"
$sFilter.Range("B3").Value = "Mid"
$sFilter.Range("A1").Value = "column"
$sFilter.Range("B1").Value = "keyword"
$sFilter.Range("C1").Value = "runmode"
$sFilter.Range("C2").Value = "Y"
$sFilter.Range("C3").Value = "Y"
$sFilter.Range("A2:A3").NumberFormat = "@"
$sFilter.Range("A2").Value = "1"
$sFilter.Range("A3").Value = "2"

# --- edit_assessment_criteria remaining cells -------------------------------
$sEdit.Range("D1").Value = "new sub assessment"
$sEdit.Range("D2").Value = "Lab"
$sEdit.Range("A1").Value = "row"
$sEdit.Range("B1").Value = "new code"
$sEdit.Range("C1").Value = "new name"
$sEdit.Range("E1").Value = "runmode"
$sEdit.Range("E2").Value = "Y"
$sEdit.Range("D3").Value = "Lab"
$sEdit.Range("E3").Value = "Y"
$sEdit.Range("A2:A3").NumberFormat = "@"
$sEdit.Range("A2").Value = "1"
$sEdit.Range("A3").Value = "2"

$sEdit.Columns.Item(2).ColumnWidth = 17.451822916666668
$sEdit.Columns.Item(3).ColumnWidth = 18.166666666666668
$sEdit.Columns.Item(4).ColumnWidth = 19.022135416666668

# --- view_assessment_criteria ------------------------------------------------
$sView.Range("A1").Value = "row"
$sView.Range("B1").Value = "runmode"
$sView.Range("B2").Value = "Y"
$sView.Range("B3").Value = "Y"
$sView.Range("A2:A3").NumberFormat = "@"
$sView.Range("A2").Value = "1"
$sView.Range("A3").Value = "2"

# ---------------------------------------------------------------------------
# 3. Selections on the new sheets (last Select() call wins the active tab)
# ---------------------------------------------------------------------------

$sAdd.Range("E4").Select() | Out-Null
$sSearch.Range("C5").Select() | Out-Null
$sView.Range("A1:B3").Select() | Out-Null
$sFilter.Range("C2:C3").Select() | Out-Null
$sEdit.Range("A2:A3").Select() | Out-Null
